$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record (Feria Lagunitas de Puerto Montt - Haba) is inserted
# as row 41, pushing the existing rows 41:52 down to 42:53.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new record's data.
$ws.Cells.Item(41, 1).Value = 4
$ws.Cells.Item(41, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(41, 3).Value = "Los Lagos"
$ws.Cells.Item(41, 4).Value = 44508
$ws.Cells.Item(41, 5).Value = 10
$ws.Cells.Item(41, 6).Value = 100112026
$ws.Cells.Item(41, 7).Value = "Haba"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 80
$ws.Cells.Item(41, 11).Value = 10000
$ws.Cells.Item(41, 12).Value = 10000
$ws.Cells.Item(41, 13).Value = 10000
$ws.Cells.Item(41, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(41, 15).Value = "Región del Maule"
$ws.Cells.Item(41, 16).Value = 400
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"
